$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.429.03'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.53%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.907.27'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.14%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.69'
$ws.Range("D5").Style = "Normal"

$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4680'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.90%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4079'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.81%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.63'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.86%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08008'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.13%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.006'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.31%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.30'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.83%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.910.54'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.931'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.117'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.71%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.04'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.35%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06597'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001025'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.54%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.70'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.45%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '29.477.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.63%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.524'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.40%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.54'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.97%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.171.61'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.78%  '

$ws.Range("E27").Value = '  -2.80%  '

$ws.Range("E28").Value = '  +0.31%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.129'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.50%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.703'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.72%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '116.85'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.67%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.070'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +9.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09483'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.04%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.419'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.574'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.57%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.372'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.49%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02254'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.71%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06077'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.02%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.347'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.52%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.172'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.61%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5862'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.23%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1834'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.62%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.465'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.84%  '

$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.10'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.301'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.25%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07726'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.96%  '

$ws.Range("E47").Value = '  +0.43%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5535'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.76%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.930'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.65%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '113.22'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.99%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.2943'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.15%  '
